$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B ("Rule" name for the 4th rule row) changes from "R40"
# to the text "1". Assigning a bare numeric-looking string via .Value
# would store it as a number (losing the shared-string/text nature), so
# we stage the text in a scratch cell that's explicitly formatted as
# Text, copy it, and paste-special *values only* into B11. That keeps
# B11's existing style/format untouched while swapping its content over
# to the literal text "1".
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()

$ws.Range("B11").PasteSpecial(-4163)

# Clean up the scratch cell and clipboard state so nothing else changes.
$scratch.Clear()
$excel.CutCopyMode = 0
